$d = $word.ActiveDocument

# --- 1) "Qualidades" paragraph -------------------------------------------
# Merge the "Qualidades" + ":" runs (both bold, identical formatting) into
# a single "Qualidades:" run.
$d.Content.Find.Execute("Qualidades:", $true, $false, $false, $false, $false, $true, 1, $false, "Qualidades:", 1) | Out-Null

# Shrink the placeholder list from 5 blanks to 3 for both "Qualidades" and
# "Defeitos" (identical pattern in both paragraphs) in one pass.
$d.Content.Find.Execute(" _, _, _, _, _.", $true, $false, $false, $false, $false, $true, 1, $false, " _, _, _.", 2) | Out-Null

# --- 2) "Nível" paragraph -------------------------------------------------
# Merge "Nível" + ":" + " " (all bold, identical formatting) into one run.
$d.Content.Find.Execute("Nível:", $true, $false, $false, $false, $false, $true, 1, $false, "Nível:", 1) | Out-Null

# The merge above also coalesces the following "1" and "." runs (they share
# identical non-bold formatting), which the original document keeps split.
# Re-split them by toggling a character property on "1" only and reverting
# it, which forces the run boundary to re-appear without altering the
# visible formatting.
$rngNivel = $d.Content
$foundNivel = $rngNivel.Find.Execute("1.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundNivel) {
    $digit = $d.Range($rngNivel.Start, $rngNivel.Start + 1)
    $digit.Bold = 1
    $digit.Bold = 0
}

# --- 3) "Aqua" paragraph --------------------------------------------------
# Merge "Aqua" + ":" + " " (all bold, identical formatting) into one run.
# "Grau 0." that follows is a differently-formatted (non-bold) run, so it
# is not affected by the merge.
$d.Content.Find.Execute("Aqua:", $true, $false, $false, $false, $false, $true, 1, $false, "Aqua:", 1) | Out-Null
